$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated raw "Nb nouveaux cas positifs" (column C) figures for several days ---
$ws.Range("C394").Value = 107
$ws.Range("C397").Value = 41
$ws.Range("C427").Value = 99
$ws.Range("C430").Value = 86
$ws.Range("C432").Value = 29
$ws.Range("C433").Value = 111
$ws.Range("C434").Value = 89
$ws.Range("C435").Value = 40

# --- Row 435 (2021-05-05) corrections ---
$ws.Range("E435").Value = 8
$ws.Range("F435").Value = 6
$ws.Range("L435").Value = 1

# --- Row 436 (2021-05-06) newly filled in with the day's figures ---
$ws.Range("C436").Value = 1
$ws.Range("E436").Value = 8
$ws.Range("F436").Value = 7
$ws.Range("G436").Value = 22
$ws.Range("L436").Value = 0
$ws.Range("M436").Value = 0

# --- View state: scroll position / active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 2
$ws.Range("Q23").Select()
